$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-22: update Price (D) and Volume(1h) (E) columns
$ws.Range("D2").Value = "25.520.77"
$ws.Range("E2").Value = "  -6.34%  "
$ws.Range("D3").Value = "1.805.05"
$ws.Range("E3").Value = "  -5.33%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'275.67"
$ws.Range("E5").Value = "  -10.03%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.5007"
$ws.Range("E7").Value = "  -6.85%  "
$ws.Range("D8").Value = "'0.3497"
$ws.Range("E8").Value = "  -8.37%  "
$ws.Range("D9").Value = "'43.81"
$ws.Range("E9").Value = "  -4.70%  "
$ws.Range("D10").Value = "'0.06663"
$ws.Range("E10").Value = "  -8.67%  "
$ws.Range("D11").Value = "'19.94"
$ws.Range("E11").Value = "  -10.52%  "
$ws.Range("D12").Value = "'0.8338"
$ws.Range("E12").Value = "  -7.95%  "
$ws.Range("D13").Value = "'0.07816"
$ws.Range("E13").Value = "  -4.69%  "
$ws.Range("D14").Value = "1.811.19"
$ws.Range("E14").Value = "  +63.18%  "
$ws.Range("D15").Value = "'5.048"
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("D16").Value = "'87.13"
$ws.Range("E16").Value = "  -9.11%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "'13.89"
$ws.Range("E18").Value = "  -6.60%  "
$ws.Range("D19").Value = "'1.002"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'0.000007918"
$ws.Range("E20").Value = "  -8.77%  "
$ws.Range("D21").Value = "25.606.52"
$ws.Range("E21").Value = "  -6.09%  "
$ws.Range("D22").Value = "'4.710"
$ws.Range("E22").Value = "  -6.73%  "

# Rows 23-51: replace Coin (B), Link (C), Price (D), Volume(1h) (E) columns
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'9.943"
$ws.Range("E23").Value = "  -7.89%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'6.051"
$ws.Range("E24").Value = "  -7.29%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "'141.33"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "'2.112"
$ws.Range("E26").Value = "  -8.50%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.655"
$ws.Range("E27").Value = "  -5.58%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'16.94"
$ws.Range("E28").Value = "  -8.04%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'108.70"
$ws.Range("E29").Value = "  -6.90%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.294"
$ws.Range("E30").Value = "  -11.18%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.204"
$ws.Range("E31").Value = "  -11.06%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08861"
$ws.Range("E32").Value = "  -3.98%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04808"
$ws.Range("E33").Value = "  -5.42%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7305"
$ws.Range("E34").Value = "  -11.86%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.122"
$ws.Range("E35").Value = "  -7.92%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.842"
$ws.Range("E36").Value = "  -5.11%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'1.001"
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'3.025"
$ws.Range("E38").Value = "  -8.79%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01856"
$ws.Range("E39").Value = "  -7.39%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5179"
$ws.Range("E40").Value = "  -12.91%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.258"
$ws.Range("E41").Value = "  -15.44%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9592"
$ws.Range("E42").Value = "  -11.08%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'112.84"
$ws.Range("E43").Value = "  -3.38%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.165"
$ws.Range("E44").Value = "  -7.69%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.035"
$ws.Range("E45").Value = "  -14.30%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4577"
$ws.Range("E47").Value = "  -11.37%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1379"
$ws.Range("E48").Value = "  -10.03%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.199"
$ws.Range("E49").Value = "  -10.09%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'35.65"
$ws.Range("E50").Value = "  -7.32%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.486"
$ws.Range("E51").Value = "  -9.68%  "
